$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.824.17"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.088.53"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +2.51%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0792"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.399.02"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.81"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.776"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.092.87"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.736.63"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.23"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.50"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("E27").Value = "  +4.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.05"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.42"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.62"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.73"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0636"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.51"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.45"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0982"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.33"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.91"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0216"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.463.88"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.08"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.27%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.47"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.283.09"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.53%  "
